$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 12): update Prediction / Error
$ws.Range("D2").Value = 0.9999995250420277
$ws.Range("E2").Value = 0.9999995250420277

# Row 3 (Control 18): Success flips to TRUE, Prediction / Error updated
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.09155899752677703
$ws.Range("E3").Value = 0.09155899752677703

# Row 4 (Control 34): Prediction / Error updated
$ws.Range("D4").Value = 0.00002062509513299652
$ws.Range("E4").Value = 0.00002062509513299652

# Row 5 (Control 42): Prediction / Error updated
$ws.Range("D5").Value = 0.0000000000000000000000000000000363018590281252
$ws.Range("E5").Value = 0.0000000000000000000000000000000363018590281252

# Row 7 (MDD 27): Prediction / Error updated
$ws.Range("D7").Value = 0.9999999157632988
$ws.Range("E7").Value = 0.00000008423670116997073

# Row 8 (MDD 47): Prediction / Error updated
$ws.Range("D8").Value = 0.00001218410908539108
$ws.Range("E8").Value = 0.9999878158909146

# Row 9 (MDD 13): Prediction / Error updated
$ws.Range("D9").Value = 0.02422699046677796
$ws.Range("E9").Value = 0.975773009533222

# Row 11 (MDD 5): Success flips to FALSE, Prediction / Error / Cross Entropy Loss updated
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.3772118311264953
$ws.Range("E11").Value = 0.6227881688735046
$ws.Range("F11").Value = 11.28780937194824
